$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.336.82"
$ws.Range("E2").Value = "'  -0.08%  "
$ws.Range("D3").Value = "'1.933.54"
$ws.Range("E3").Value = "'  -0.24%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "'  -0.13%  "
$ws.Range("D5").Value = "'0.7473"
$ws.Range("E5").Value = "'  +3.11%  "
$ws.Range("D6").Value = "'248.95"
$ws.Range("E6").Value = "'  -0.59%  "
$ws.Range("E7").Value = "'  -0.05%  "
$ws.Range("D8").Value = "'28.20"
$ws.Range("E8").Value = "'  -0.66%  "
$ws.Range("D9").Value = "'0.3208"
$ws.Range("E9").Value = "'  -3.80%  "
$ws.Range("D10").Value = "'0.07115"
$ws.Range("E10").Value = "'  -1.77%  "
$ws.Range("D11").Value = "'0.7885"
$ws.Range("E11").Value = "'  -2.68%  "
$ws.Range("D12").Value = "'0.08003"
$ws.Range("E12").Value = "'  -1.19%  "
$ws.Range("D13").Value = "'1.933.54"
$ws.Range("E13").Value = "'  -0.25%  "
$ws.Range("D14").Value = "'5.386"
$ws.Range("E14").Value = "'  -1.58%  "
$ws.Range("D15").Value = "'94.50"
$ws.Range("E15").Value = "'  -0.02%  "
$ws.Range("D16").Value = "'14.63"
$ws.Range("E16").Value = "'  -2.55%  "
$ws.Range("D17").Value = "'30.343.72"
$ws.Range("E17").Value = "'  -0.07%  "
$ws.Range("D18").Value = "'252.73"
$ws.Range("E18").Value = "'  +1.18%  "
$ws.Range("D19").Value = "'0.000008040"
$ws.Range("E19").Value = "'  -2.35%  "
$ws.Range("D20").Value = "'5.788"
$ws.Range("E20").Value = "'  -2.24%  "
$ws.Range("D21").Value = "'2.188.25"
$ws.Range("E21").Value = "'  -0.08%  "
$ws.Range("E22").Value = "'  +0.00%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "'  -0.16%  "
$ws.Range("E24").Value = "'  -1.66%  "
$ws.Range("D25").Value = "'9.573"
$ws.Range("E25").Value = "'  -1.94%  "
$ws.Range("D26").Value = "'164.57"
$ws.Range("E26").Value = "'  +0.71%  "
$ws.Range("D27").Value = "'2.332"
$ws.Range("E27").Value = "'  -2.42%  "
$ws.Range("D28").Value = "'19.12"
$ws.Range("E28").Value = "'  -0.96%  "
$ws.Range("D29").Value = "'0.1322"
$ws.Range("E29").Value = "'  -0.47%  "
$ws.Range("E30").Value = "'  +0.78%  "
$ws.Range("D31").Value = "'1.532"
$ws.Range("E31").Value = "'  -2.46%  "
$ws.Range("D32").Value = "'4.443"
$ws.Range("E32").Value = "'  -0.02%  "
$ws.Range("D33").Value = "'4.153"
$ws.Range("E33").Value = "'  -1.10%  "
$ws.Range("D34").Value = "'0.05132"
$ws.Range("E34").Value = "'  -1.26%  "
$ws.Range("D35").Value = "'1.285"
$ws.Range("E35").Value = "'  -0.49%  "
$ws.Range("D36").Value = "'0.7497"
$ws.Range("E36").Value = "'  -0.31%  "
$ws.Range("D37").Value = "'2.771"
$ws.Range("E37").Value = "'  +0.73%  "
$ws.Range("D38").Value = "'0.01968"
$ws.Range("E38").Value = "'  -0.57%  "
$ws.Range("D39").Value = "'2.805"
$ws.Range("E39").Value = "'  -1.22%  "
$ws.Range("D40").Value = "'78.21"
$ws.Range("E40").Value = "'  -3.12%  "
$ws.Range("D41").Value = "'6.412"
$ws.Range("E41").Value = "'  -1.60%  "
$ws.Range("D42").Value = "'0.4509"
$ws.Range("E42").Value = "'  -0.74%  "
$ws.Range("D43").Value = "'1.992"
$ws.Range("E43").Value = "'  -2.34%  "
$ws.Range("D44").Value = "'0.8439"
$ws.Range("E44").Value = "'  -0.53%  "
$ws.Range("E45").Value = "'  +0.00%  "
$ws.Range("D46").Value = "'102.54"
$ws.Range("E46").Value = "'  +0.22%  "
$ws.Range("D47").Value = "'9.841"
$ws.Range("E47").Value = "'  +0.28%  "
$ws.Range("D48").Value = "'7.543"
$ws.Range("E48").Value = "'  +1.21%  "
$ws.Range("D49").Value = "'37.49"
$ws.Range("E49").Value = "'  +1.64%  "
$ws.Range("D50").Value = "'984.13"
$ws.Range("E50").Value = "'  +12.29%  "
$ws.Range("D51").Value = "'0.1192"
$ws.Range("E51").Value = "'  +4.59%  "

# Reset styles on the touched range so Excel drops the auto-added
# text-format style index that gets created when assigning values
# that look numeric (keeps cells without an explicit style, like source).
$ws.Range("D2:E51").Style = "Normal"
